$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (17) below the existing data (row 16), copying the
# formatting of row 16 (so the date cell keeps the short-date style)
# and then filling in the new values.
$ws.Range("A16:B16").Copy()
$ws.Range("A17:B17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(17, 1).Value = 45981
$ws.Cells.Item(17, 2).Value = 1

# Match the author's final selection/active cell in the saved file.
$ws.Range("E15").Select()
